# Apply the data refresh for Mon_TI_2024_02_11:
#  - updates various stat columns (5M/15M/Saison averages, age-bucket counts)
#  - adds the transposition (shift) of the M-1..M-5 columns for rows 2-4
#  - updates the M-1_vs / M-1_score / M-2_vs / M-2_score columns
#  - updates the "vs or @" / Delta / delta_B2B / nombre_de_B2B columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (Chet Holmgren)
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = 35.2

$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 6

$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 44
$ws.Range("O2").Value = 32
$ws.Range("P2").Value = 30
$ws.Range("Q2").Value = 40

$ws.Range("T2").Value = "vs"
$ws.Range("U2").Value = 32

$ws.Range("AH2").Value = "vs"
$ws.Range("AI2").Value = 1.1
$ws.Range("AJ2").Value = -7.4
$ws.Range("AK2").Value = 8

# ---------------------------------------------------------------------------
# Row 3 (Jaylen Brown)
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = 30.6
$ws.Range("F3").Value = 32.5
$ws.Range("G3").Value = 30.6

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3

$ws.Range("N3").Value = 28
$ws.Range("O3").Value = "-"
$ws.Range("P3").Value = 14
$ws.Range("Q3").Value = 41

$ws.Range("U3").Value = 30
$ws.Range("V3").Value = "vs"
$ws.Range("W3").Value = 31

$ws.Range("AH3").Value = "@"
$ws.Range("AI3").Value = -1.8

# ---------------------------------------------------------------------------
# Row 4 (De'Aaron Fox)
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = 30.2
$ws.Range("F4").Value = 30.7
$ws.Range("G4").Value = 35.4

$ws.Range("M4").Value = 28
$ws.Range("N4").Value = 15
$ws.Range("O4").Value = 17
$ws.Range("P4").Value = 60
$ws.Range("Q4").Value = 31

$ws.Range("U4").Value = 55

$ws.Range("AH4").Value = "@"
$ws.Range("AI4").Value = 0.6
